# Remove the "is_locked" and "is_enabled" columns (D and E) from the dept
# import template. Deleting the entire columns shifts the subsequent
# "order_by" and "rem" columns (F, G) left into D, E, and the now-unused
# shared strings for is_locked_lbl / is_enabled_lbl are dropped on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1:E1").EntireColumn.Delete()
